$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8848016262054443
$ws.Range("B1").Value = 2.648030519485474
$ws.Range("C1").Value = 2.445246696472168
$ws.Range("D1").Value = 1.836178064346313
$ws.Range("E1").Value = 1.345170259475708
